## Rename read_ini.ini to path.ini
## -------------------------------------------------------------
## Functional edit: append a new test case (tc007, "create interest
## group success") plus a block of placeholder rows (tc008..tc029)
## to the pbf_case test-case sheet, matching the existing table's
## look & feel (copy formatting from the last populated row), and
## touch up the window/view metadata that Excel records when a user
## does this interactively.
## -------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. new fully-populated case row (row 8) ------------------------------
$ws.Range("A8").Value = "tc007"
$ws.Range("B8").Value = "创建兴趣小组成功"
$ws.Range("C8").Value = "create_interest_group_success_p"
$ws.Range("D8").Value = "create_interest_group_success_e"

# --- 2. placeholder case rows (rows 9-30): case id filled in, params / -----
#        expect columns left blank but still bordered like the rest of the
#        table -----------------------------------------------------------
$ws.Range("A9").Value  = "tc008"
$ws.Range("A10").Value = "tc009"
$ws.Range("A11").Value = "tc010"
$ws.Range("A12").Value = "tc011"
$ws.Range("A13").Value = "tc012"
$ws.Range("A14").Value = "tc013"
$ws.Range("A15").Value = "tc014"
$ws.Range("A16").Value = "tc015"
$ws.Range("A17").Value = "tc016"
$ws.Range("A18").Value = "tc017"
$ws.Range("A19").Value = "tc018"
$ws.Range("A20").Value = "tc019"
$ws.Range("A21").Value = "tc020"
$ws.Range("A22").Value = "tc021"
$ws.Range("A23").Value = "tc022"
$ws.Range("A24").Value = "tc023"
$ws.Range("A25").Value = "tc024"
$ws.Range("A26").Value = "tc025"
$ws.Range("A27").Value = "tc026"
$ws.Range("A28").Value = "tc027"
$ws.Range("A29").Value = "tc028"
$ws.Range("A30").Value = "tc029"

# --- 3. carry the existing table formatting down onto the new rows -------
# column A uses the case_id style (font w/ explicit black + left/top align)
$ws.Range("A7").Copy()
$ws.Range("A8:A30").PasteSpecial(-4122)

# columns B:D use the plain bordered style; row 8 has real content, rows
# 9-30 stay empty (placeholders) but keep the same bordered look
$ws.Range("B7:D7").Copy()
$ws.Range("B8:D8").PasteSpecial(-4122)
$ws.Range("B7:D7").Copy()
$ws.Range("B9:D30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. widen the params/expect columns now that they hold long strings --
$ws.Columns.Item(3).ColumnWidth = 28.8
$ws.Columns.Item(4).ColumnWidth = 28.65

# --- 5. drop the stale row-outline level left over from the old sheet ----
#        (column outline level 3 is kept)
$ws.Outline.ShowLevels(0, 3)

# --- 6. move the selection down to where the user ended up editing -------
$ws.Range("B24").Select() | Out-Null

# --- 7. restore the window size recorded by the workbook view ------------
$win = $wb.Windows.Item(1)
$win.Width = 26310
$win.Height = 9480
